# The deck's slide master theme (ppt/theme/theme1.xml, "Integral") is
# switched to the stock "Office Theme" colour scheme - i.e. the 12-slot
# theme colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) that
# backs every slide is replaced with the default Office palette. The
# theme's font scheme / format scheme are untouched (they were already
# identical between the two themes present in the deck).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme colour scheme, in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Item(1).RGB  = 0x000000   # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # dk2      44546A (BGR-packed)
$tcs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6 (BGR-packed)
$tcs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5 (BGR-packed)
$tcs.Item(6).RGB  = 0x317DED   # accent2  ED7D31 (BGR-packed)
$tcs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5 (BGR-packed)
$tcs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000 (BGR-packed)
$tcs.Item(9).RGB  = 0xC47244   # accent5  4472C4 (BGR-packed)
$tcs.Item(10).RGB = 0x47AD70   # accent6  70AD47 (BGR-packed)
$tcs.Item(11).RGB = 0xC16305   # hlink    0563C1 (BGR-packed)
$tcs.Item(12).RGB = 0x724F95   # folHlink 954F72 (BGR-packed)
